$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.232145190238953
$ws.Range("B1").Value = 5.213141918182373
$ws.Range("C1").Value = 1.312560677528381
$ws.Range("D1").Value = 0.8791686296463013
$ws.Range("E1").Value = 0.7065007090568542
